$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.588.52"
$ws.Range("E2").Value = "  +4.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.744.18"
$ws.Range("E3").Value = "  +4.54%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.14"
$ws.Range("E5").Value = "  +3.85%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4820"
$ws.Range("E7").Value = "  +1.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2693"
$ws.Range("E8").Value = "  +3.49%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06262"
$ws.Range("E9").Value = "  +1.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.743.86"
$ws.Range("E10").Value = "  +4.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07135"
$ws.Range("E11").Value = "  +1.73%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.89"
$ws.Range("E12").Value = "  +7.58%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6252"
$ws.Range("E13").Value = "  +6.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.528"
$ws.Range("E14").Value = "  +3.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.48"
$ws.Range("E15").Value = "  +2.90%  "

$ws.Range("E16").Value = "  -0.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.586.82"
$ws.Range("E17").Value = "  +4.08%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.00%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006919"
$ws.Range("E19").Value = "  +2.78%  "

$ws.Range("E20").Value = "  +2.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.968.56"
$ws.Range("E21").Value = "  +4.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.627"
$ws.Range("E22").Value = "  +4.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.854"
$ws.Range("E23").Value = "  +0.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.369"
$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.51"
$ws.Range("E25").Value = "  -0.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.39"
$ws.Range("E26").Value = "  +2.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.820"
$ws.Range("E27").Value = "  +6.21%  "

$ws.Range("E28").Value = "  +3.75%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.89"
$ws.Range("E29").Value = "  +2.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.015"
$ws.Range("E30").Value = "  +0.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.741"
$ws.Range("E31").Value = "  +3.29%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.07889"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04605"
$ws.Range("E33").Value = "  +6.55%  "

$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9996"
$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.617"
$ws.Range("E35").Value = "  -0.30%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6407"
$ws.Range("E36").Value = "  +5.94%  "

$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9994"
$ws.Range("E37").Value = "  +4.86%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9312"
$ws.Range("E38").Value = "  -2.68%  "

$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "113.79"
$ws.Range("E39").Value = "  +14.85%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.998"
$ws.Range("E40").Value = "  +8.29%  "

$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.438"
$ws.Range("E41").Value = "  -3.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.813"
$ws.Range("E42").Value = "  +19.14%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.24%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01512"
$ws.Range("E44").Value = "  +2.74%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3924"
$ws.Range("E45").Value = "  +4.80%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1224"
$ws.Range("E46").Value = "  +9.83%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.784"
$ws.Range("E47").Value = "  +9.45%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05339"
$ws.Range("E48").Value = "  +1.40%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.967"
$ws.Range("E49").Value = "  +7.44%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.81"
$ws.Range("E50").Value = "  +3.00%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.260"
$ws.Range("E51").Value = "  +4.88%  "
